# Update latest output (run 34)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E2").Value = 579.8138287500001
$schedule.Range("F2").Value = 12.78249181547619
$schedule.Range("E3").Value = 378.08765475
$schedule.Range("F3").Value = 25.0057972718254

# --- Sheet "Detailed" ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B17").Value = 56.98
$detailed.Range("B18").Value = 57.06003
$detailed.Range("B19").Value = 36.06
$detailed.Range("C19").Value = "historical"
$detailed.Range("C20").Value = "historical"
$detailed.Range("B21").Value = 36.06
$detailed.Range("B22").Value = -5.01
$detailed.Range("B23").Value = 22.07
$detailed.Range("B24").Value = 36.06
$detailed.Range("B25").Value = -16.93619
$detailed.Range("B26").Value = -21.63841
$detailed.Range("B27").Value = -21.91848
$detailed.Range("B28").Value = -17.11981
$detailed.Range("B29").Value = -18.87147
$detailed.Range("B30").Value = -23.5
$detailed.Range("B31").Value = -19.98
$detailed.Range("B32").Value = -16.78441
$detailed.Range("B33").Value = -16.79274
$detailed.Range("B34").Value = 19.11478
$detailed.Range("B35").Value = 7.901
$detailed.Range("B36").Value = -9.30804
$detailed.Range("B37").Value = -8.08661
$detailed.Range("B38").Value = -8.467040000000001
$detailed.Range("B39").Value = -3.07664
$detailed.Range("B40").Value = -3.07313
$detailed.Range("B41").Value = 29.85322
$detailed.Range("B45").Value = 64.8901
$detailed.Range("B48").Value = 57.03042
